$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 216.58333
$ws.Range("I6").Value = 216.58333
$ws.Range("K6").Value = 649.74999
$ws.Range("M6").Value = -537.74999

$ws.Range("H8").Value = 14.833333
$ws.Range("I8").Value = 16
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = 48
$ws.Range("L8").Value = 6
$ws.Range("M8").Value = 91
$ws.Range("N8").Value = -284

$ws.Range("H9").Value = 109.84615
$ws.Range("I9").Value = 39.444443
$ws.Range("J9").Value = 268.25
$ws.Range("K9").Value = 39.444443
$ws.Range("L9").Value = 268.25
$ws.Range("M9").Value = 129.555557
$ws.Range("N9").Value = -606.25

$ws.Range("H15").Value = 90.2
$ws.Range("I15").Value = 90.2
$ws.Range("K15").Value = 270.6
$ws.Range("M15").Value = -101.6

$ws.Range("H21").Value = 2409.4
$ws.Range("I21").Value = 2409.4
$ws.Range("K21").Value = 2409.4
$ws.Range("M21").Value = -1941.4

$ws.Range("H23").Value = 2409.4
$ws.Range("I23").Value = 2409.4
$ws.Range("K23").Value = 2409.4
$ws.Range("M23").Value = -2175.4

$ws.Range("H75").Value = 9000
$ws.Range("I75").Value = 8500
$ws.Range("J75").Value = 9500
$ws.Range("K75").Value = 8500
$ws.Range("L75").Value = 9500
$ws.Range("M75").Value = -7564
$ws.Range("N75").Value = -11372

$ws.Range("H78").Value = 9000
$ws.Range("I78").Value = 8500
$ws.Range("J78").Value = 9500
$ws.Range("K78").Value = 25500
$ws.Range("L78").Value = 28500
$ws.Range("M78").Value = -20820
$ws.Range("N78").Value = -37860

$ws.Range("H88").Value = 3357.4285
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3357.4285
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3357.4285
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -4169.4285

$ws.Range("H91").Value = 3357.4285
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3357.4285
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3357.4285
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -6165.4285

$ws.Range("H113").Value = 30999.4
$ws.Range("I113").Value = 36499.25
$ws.Range("J113").Value = 9000
$ws.Range("K113").Value = 36499.25
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -33245.25
$ws.Range("N113").Value = -15508

$ws.Range("H137").Value = 1498.7273
$ws.Range("J137").Value = 1374.5
$ws.Range("L137").Value = 4123.5
$ws.Range("N137").Value = -9223.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H74").Value = 1293.125
$ws.Range("I74").Value = 1206.4286
$ws.Range("J74").Value = 1900
$ws.Range("K74").Value = 1206.4286
$ws.Range("L74").Value = 1900
$ws.Range("M74").Value = -332.4286
$ws.Range("N74").Value = -3648

$ws.Range("H77").Value = 1293.125
$ws.Range("I77").Value = 1206.4286
$ws.Range("J77").Value = 1900
$ws.Range("K77").Value = 6032.143
$ws.Range("L77").Value = 9500
$ws.Range("M77").Value = -1664.143
$ws.Range("N77").Value = -18236

$ws.Range("H122").Value = 1865.6666
$ws.Range("I122").Value = 1865.6666
$ws.Range("K122").Value = 5596.9998
$ws.Range("M122").Value = -3146.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1999.5
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 1999.5
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 1999.5
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -2493.5

$ws.Range("H99").Value = 2741
$ws.Range("J99").Value = 3372
$ws.Range("L99").Value = 3372
$ws.Range("N99").Value = -6368

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1690.4
$ws.Range("I7").Value = 1095.7
$ws.Range("J7").Value = 2879.8
$ws.Range("K7").Value = 1095.7
$ws.Range("L7").Value = 2879.8
$ws.Range("M7").Value = -982.7
$ws.Range("N7").Value = -3105.8

$ws.Range("H31").Value = 4623
$ws.Range("I31").Value = 3872.7
$ws.Range("K31").Value = 3872.7
$ws.Range("M31").Value = -3577.7

$ws.Range("H34").Value = 4623
$ws.Range("I34").Value = 3872.7
$ws.Range("K34").Value = 3872.7
$ws.Range("M34").Value = -3670.7

$ws.Range("H99").Value = 2871.7144
$ws.Range("J99").Value = 3664.6667
$ws.Range("L99").Value = 3664.6667
$ws.Range("N99").Value = -6660.6667

$ws.Range("H126").Value = 2871.7144
$ws.Range("J126").Value = 3664.6667
$ws.Range("L126").Value = 10994.0001
$ws.Range("N126").Value = -15934.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 412.9091
$ws.Range("I7").Value = 612.4286
$ws.Range("J7").Value = 63.75
$ws.Range("K7").Value = 1837.2858
$ws.Range("L7").Value = 191.25
$ws.Range("M7").Value = -1725.2858
$ws.Range("N7").Value = -415.25

$ws.Range("H17").Value = 849.5
$ws.Range("J17").Value = 849.5
$ws.Range("L17").Value = 2548.5
$ws.Range("N17").Value = -2886.5

$ws.Range("H119").Value = 832
$ws.Range("I119").Value = 832
$ws.Range("K119").Value = 2496
$ws.Range("M119").Value = 2342

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2017.8823
$ws.Range("I97").Value = 1996.091
$ws.Range("J97").Value = 2057.8333
$ws.Range("K97").Value = 1996.091
$ws.Range("L97").Value = 2057.8333
$ws.Range("M97").Value = -1500.091
$ws.Range("N97").Value = -3049.8333

$ws.Range("H132").Value = 3968.3809
$ws.Range("I132").Value = 3866.8
$ws.Range("K132").Value = 11600.4
$ws.Range("M132").Value = -9070.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 7500
$ws.Range("I14").Value = 10000
$ws.Range("K14").Value = 10000
$ws.Range("M14").Value = -9828

$ws.Range("H40").Value = 4243
$ws.Range("I40").Value = 4243
$ws.Range("K40").Value = 4243
$ws.Range("M40").Value = -4107

$ws.Range("H55").Value = 933.3333
$ws.Range("I55").Value = 925
$ws.Range("K55").Value = 925
$ws.Range("M55").Value = -752

$ws.Range("H82").Value = 976.2353000000001
$ws.Range("J82").Value = 1170.7142
$ws.Range("L82").Value = 1170.7142
$ws.Range("N82").Value = -1892.7142

$ws.Range("H85").Value = 976.2353000000001
$ws.Range("J85").Value = 1170.7142
$ws.Range("L85").Value = 1170.7142
$ws.Range("N85").Value = -3666.7142

$ws.Range("H93").Value = 394.66666
$ws.Range("I93").Value = 400
$ws.Range("K93").Value = 400
$ws.Range("M93").Value = 848

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H136").Value = 4666.3335
$ws.Range("I136").Value = 6000
$ws.Range("J136").Value = 1999
$ws.Range("K136").Value = 18000
$ws.Range("L136").Value = 5997
$ws.Range("M136").Value = -15450
$ws.Range("N136").Value = -11097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 4575
$ws.Range("J6").Value = 6750
$ws.Range("L6").Value = 6750
$ws.Range("N6").Value = -6980

$ws.Range("H11").Value = 3352084
$ws.Range("I11").Value = 10001002
$ws.Range("J11").Value = 27625
$ws.Range("K11").Value = 10001002
$ws.Range("L11").Value = 27625
$ws.Range("M11").Value = -10000860
$ws.Range("N11").Value = -27909

$ws.Range("H20").Value = 28225
$ws.Range("I20").Value = 28225
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 28225
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -27985
$ws.Range("N20").ClearContents()

$ws.Range("H100").Value = 3873364.8
$ws.Range("I100").Value = 5809047.5
$ws.Range("J100").Value = 1999.1666
$ws.Range("K100").Value = 11618095
$ws.Range("L100").Value = 3998.3332
$ws.Range("M100").Value = -11617554
$ws.Range("N100").Value = -5080.3332
